# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.136.25"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").Value = "3.138.27"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.14"
$ws.Range("E5").Value = "  +3.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.12"
$ws.Range("E6").Value = "  +2.90%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").Value = "  +11.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.32"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.425"
$ws.Range("E11").Value = "  +6.56%  "
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("D13").Value = "3.678.35"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.99"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("E15").Value = "  +5.39%  "
$ws.Range("D16").Value = "58.235.43"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.25"
$ws.Range("E17").Value = "  +6.77%  "
$ws.Range("D18").Value = "3.126.85"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("E19").Value = "  +4.51%  "
$ws.Range("E20").Value = "  +4.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.26"
$ws.Range("E21").Value = "  +7.27%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.75"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.45"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("E25").Value = "  +3.78%  "
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.12"
$ws.Range("E28").Value = "  +13.03%  "
$ws.Range("D29").Value = "0.0₃0885"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("E31").Value = "  +6.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.79"
$ws.Range("E32").Value = "  +4.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.16"
$ws.Range("E33").Value = "  +7.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.18"
$ws.Range("E34").Value = "  +4.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.26"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.27"
$ws.Range("E36").Value = "  +4.86%  "
$ws.Range("E37").Value = "  +12.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.54"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.69"
$ws.Range("E39").Value = "  +7.31%  "
$ws.Range("D40").Value = "2.647.94"
$ws.Range("E40").Value = "  +10.51%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0679"
$ws.Range("E41").Value = "  +3.85%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.26"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.69"
$ws.Range("E43").Value = "  +5.66%  "
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("E45").Value = "  +4.94%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.102"
$ws.Range("E47").Value = "  +11.89%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.24"
$ws.Range("E48").Value = "  +4.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.976"
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.31"
$ws.Range("E50").Value = "  +3.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.748"
$ws.Range("E51").Value = "  -0.11%  "
